$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 91) — DGS's 2021/10/08 report.
# Column A holds the date as text (shared string), matching the existing
# rows above it, so force a text format before assigning the value to
# prevent Excel from auto-converting the "yyyy/mm/dd"-looking text into a
# real date serial number, then restore the original date display format.
$row = 91
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2021/10/08"
$ws.Cells.Item($row, 1).NumberFormat = "yyyy/mm/dd"

$ws.Cells.Item($row, 2).Value = 86.5
$ws.Cells.Item($row, 3).Value = 86.7
$ws.Cells.Item($row, 4).Value = 0.92
$ws.Cells.Item($row, 5).Value = 0.91

# Move the active selection down to the next empty row, as Excel would
# after a user types the row of data and presses Enter.
$ws.Range("A92").Select() | Out-Null
